$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Shared-string text change: "Ready for handoff" -> "In Translation"
#    This status string is used on:
#      - Overview!E2, F2, E3, F3   (per-locale status columns)
#      - zh-cn!C2, C3              (Status column)
#      - de-de!C2, C3              (Status column)
#    Updating every cell that currently shows "Ready for handoff" collapses
#    them back onto a single shared string once saved.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$statusOld = "Ready for handoff"
$statusNew = "In Translation"

foreach ($addr in @("E2", "F2", "E3", "F3")) {
    $cell = $overview.Range($addr)
    if ($cell.Value2 -eq $statusOld) {
        $cell.Value = $statusNew
    }
}

foreach ($addr in @("C2", "C3")) {
    $cell = $zhcn.Range($addr)
    if ($cell.Value2 -eq $statusOld) {
        $cell.Value = $statusNew
    }
}

foreach ($addr in @("C2", "C3")) {
    $cell = $dede.Range($addr)
    if ($cell.Value2 -eq $statusOld) {
        $cell.Value = $statusNew
    }
}

# ---------------------------------------------------------------------------
# 2) Narrow the per-locale "Status" columns.
#    Overview columns E (zh-cn) and F (de-de), plus the "Status" column (C)
#    on each per-locale sheet, shrink from ~17.22 to ~13.41 characters wide.
# ---------------------------------------------------------------------------
$newStatusWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newStatusWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusWidth

$zhcn.Columns.Item(3).ColumnWidth = $newStatusWidth
$dede.Columns.Item(3).ColumnWidth = $newStatusWidth
